$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NPC")

# Insert a new column before column E (GPT_MODEL). This shifts
# GPT_MODEL/PORT/API/attributes one column to the right (E->F, F->G, G->H, H->I)
# and creates the new "conversation_example" column at E.
$ws.Columns.Item(5).Insert()

$ws.Cells.Item(1, 5).Value2 = "conversation_example"
$ws.Cells.Item(4, 5).Value2 = "赞美光明！教宗也许是多虑了。任何异变都不可能对教廷产生威胁。即使有。哼哼～也会被我一如既往的‘抹除’掉。"
$ws.Cells.Item(5, 5).Value2 = "嘶嘶嘶。。。。 去吧! 我的子嗣，去吧！为我带来食物！为我监视一切。如果有人背叛我，就让他们知道背叛的代价！嘶嘶嘶！`n任何档胆敢进入腐臭地窖的活物，都要被我吃掉！嘶嘶嘶。。。。哦，对了。除了我的小玩应——好运气先生，这个小家伙还是有点用的。如果它叫醒了我，也许是出了什么事。`n不过，没关系，嘶嘶嘶，任何的胆敢忤逆我的，都将被我撕成碎片！嘶嘶嘶！！！！"
$ws.Cells.Item(6, 5).Value2 = "吱吱。。。。格雷和他的该死的狗——摩尔。我讨厌他们，真心讨厌他们！吱吱！！鼠王没吃掉他们，还不是他们定期给鼠王送‘食物’？呵呵。吱吱，无耻的东西！`n走着瞧吧！吱吱。。。。早晚有一天他们没用了，鼠王就会吃掉他们，希望到时候他们别瘦的没肉可吃。吱吱。。。。我要时刻盯着他们，因为如果出了问题，鼠王连我也不会放过！吱吱！！`n总之，出了任何异常我都会去腐臭地窖叫醒鼠王！吱吱！！"
$ws.Cells.Item(7, 5).Value2 = "（阴郁的有气无力的声音）哦，我知道了。。。嗯？这种诅咒的命运我已经习惯了。最近没有给鼠王送食物，他似乎不太高兴。`n让我想想。该怎么办呢？要不，把其他的活人骗进腐臭地窖好了(狡诈的坏笑)，嘻嘻嘻嘻。。。那一定很有趣（恶毒的笑声）。。。总之我是不会去的，因为进入腐臭地窖的一定会被鼠王吃掉。"
$ws.Cells.Item(8, 5).Value2 = "（狗叫声）汪汪汪。。。。（狗叫声）汪汪汪。。。。（狗叫声）汪汪汪。。。。"
$ws.Cells.Item(9, 5).Value2 = "我知道了。"

# Column width adjustments observed in the target workbook (xlsx <col width="...">
# is ColumnWidth + 5/6; subtract the offset so the saved XML lands on the exact
# integer widths used by the author).
$offset = 5 / 6
$ws.Columns.Item(4).ColumnWidth = 419 - $offset
$ws.Columns.Item(5).ColumnWidth = 70 - $offset
$ws.Columns.Item(6).ColumnWidth = 49 - $offset
$ws.Columns.Item(7).ColumnWidth = 17 - $offset
$ws.Columns.Item(8).ColumnWidth = 23 - $offset
$ws.Columns.Item(21).ColumnWidth = 17 - $offset

Write-Host "Done editing NPC sheet: inserted conversation_example column."
